$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29 - this shifts existing rows 29-134 down to 30-135,
# carrying their formatting (including the date style on column D) along with them.
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44525
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112021
$ws.Range("G29").Value = "Ají"
$ws.Range("H29").Value = "Americana (o)"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 18000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 18000
$ws.Range("N29").Value = "$/caja 15 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 1200
$ws.Range("Q29").Value = 15
$ws.Range("R29").Value = "Hortaliza"
